$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2606.5386
$ws.Range("I33").Value = 1563.7
$ws.Range("J33").Value = 6082.6665
$ws.Range("K33").Value = 1563.7
$ws.Range("L33").Value = 6082.6665
$ws.Range("M33").Value = -1334.7
$ws.Range("N33").Value = -6540.6665
$ws.Range("H49").Value = 732
$ws.Range("I49").Value = 732
$ws.Range("K49").Value = 2196
$ws.Range("M49").Value = -2060
$ws.Range("H70").Value = 2300
$ws.Range("I70").Value = 2246.6667
$ws.Range("J70").Value = 2320
$ws.Range("K70").Value = 6740.000100000001
$ws.Range("L70").Value = 6960
$ws.Range("M70").Value = -6470.000100000001
$ws.Range("N70").Value = -7500
$ws.Range("H73").Value = 2300
$ws.Range("I73").Value = 2246.6667
$ws.Range("J73").Value = 2320
$ws.Range("K73").Value = 6740.000100000001
$ws.Range("L73").Value = 6960
$ws.Range("M73").Value = -5804.000100000001
$ws.Range("N73").Value = -8832
$ws.Range("H138").Value = 4074
$ws.Range("I138").Value = 3181.7273
$ws.Range("J138").Value = 13889
$ws.Range("K138").Value = 9545.1819
$ws.Range("L138").Value = 41667
$ws.Range("M138").Value = -4405.1819
$ws.Range("N138").Value = -51947

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1790.9445
$ws.Range("I45").Value = 1327.4375
$ws.Range("J45").Value = 5499
$ws.Range("K45").Value = 1327.4375
$ws.Range("L45").Value = 5499
$ws.Range("M45").Value = -950.4375
$ws.Range("N45").Value = -6253
$ws.Range("H59").Value = 65000
$ws.Range("J59").Value = 65000
$ws.Range("L59").Value = 65000
$ws.Range("N59").Value = -66608
$ws.Range("H76").Value = 159843.62
$ws.Range("I76").Value = 80000
$ws.Range("J76").Value = 171249.86
$ws.Range("K76").Value = 80000
$ws.Range("L76").Value = 171249.86
$ws.Range("M76").Value = -79662
$ws.Range("N76").Value = -171925.86
$ws.Range("H79").Value = 159843.62
$ws.Range("I79").Value = 80000
$ws.Range("J79").Value = 171249.86
$ws.Range("K79").Value = 80000
$ws.Range("L79").Value = 171249.86
$ws.Range("M79").Value = -78830
$ws.Range("N79").Value = -173589.86
$ws.Range("H122").Value = 3364.7917
$ws.Range("I122").Value = 2599.0557
$ws.Range("J122").Value = 5662
$ws.Range("K122").Value = 7797.1671
$ws.Range("L122").Value = 16986
$ws.Range("M122").Value = -5347.1671
$ws.Range("N122").Value = -21886

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1800.3889
$ws.Range("I99").Value = 1733.6364
$ws.Range("J99").Value = 1905.2858
$ws.Range("K99").Value = 1733.6364
$ws.Range("L99").Value = 1905.2858
$ws.Range("M99").Value = -235.6364000000001
$ws.Range("N99").Value = -4901.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5269.2104
$ws.Range("I62").Value = 4673.933
$ws.Range("J62").Value = 7501.5
$ws.Range("K62").Value = 4673.933
$ws.Range("L62").Value = 7501.5
$ws.Range("M62").Value = -4049.933
$ws.Range("N62").Value = -8749.5
$ws.Range("H65").Value = 5269.2104
$ws.Range("I65").Value = 4673.933
$ws.Range("J65").Value = 7501.5
$ws.Range("K65").Value = 23369.665
$ws.Range("L65").Value = 37507.5
$ws.Range("M65").Value = -20249.665
$ws.Range("N65").Value = -43747.5
$ws.Range("H99").Value = 3697.95
$ws.Range("I99").Value = 1783
$ws.Range("K99").Value = 1783
$ws.Range("M99").Value = -285
$ws.Range("H122").Value = 1964.1428
$ws.Range("I122").Value = 2031.5
$ws.Range("J122").Value = 1560
$ws.Range("K122").Value = 6094.5
$ws.Range("L122").Value = 4680
$ws.Range("M122").Value = -3644.5
$ws.Range("N122").Value = -9580
$ws.Range("H126").Value = 3697.95
$ws.Range("I126").Value = 1783
$ws.Range("K126").Value = 5349
$ws.Range("M126").Value = -2879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 79.25
$ws.Range("I2").Value = 109.42857
$ws.Range("J2").Value = 37
$ws.Range("K2").Value = 656.57142
$ws.Range("L2").Value = 222
$ws.Range("M2").Value = -543.57142
$ws.Range("N2").Value = -448
$ws.Range("H68").Value = 1894.1428
$ws.Range("I68").Value = 2899.5
$ws.Range("J68").Value = 1764.4193
$ws.Range("K68").Value = 8698.5
$ws.Range("L68").Value = 5293.257900000001
$ws.Range("M68").Value = -7887.5
$ws.Range("N68").Value = -6915.257900000001
$ws.Range("H70").Value = 14779.8
$ws.Range("I70").Value = 11299.667
$ws.Range("K70").Value = 33899.001
$ws.Range("M70").Value = -33584.001
$ws.Range("H71").Value = 1894.1428
$ws.Range("I71").Value = 2899.5
$ws.Range("J71").Value = 1764.4193
$ws.Range("K71").Value = 26095.5
$ws.Range("L71").Value = 15879.7737
$ws.Range("M71").Value = -22039.5
$ws.Range("N71").Value = -23991.7737
$ws.Range("H73").Value = 14779.8
$ws.Range("I73").Value = 11299.667
$ws.Range("K73").Value = 33899.001
$ws.Range("M73").Value = -32807.001
$ws.Range("H113").Value = 949.9091
$ws.Range("J113").Value = 994.9
$ws.Range("L113").Value = 2984.7
$ws.Range("N113").Value = -7324.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 165571.28
$ws.Range("J64").Value = 165571.28
$ws.Range("L64").Value = 165571.28
$ws.Range("N64").Value = -166067.28
$ws.Range("H67").Value = 165571.28
$ws.Range("J67").Value = 165571.28
$ws.Range("L67").Value = 165571.28
$ws.Range("N67").Value = -167287.28
$ws.Range("H69").Value = 199860
$ws.Range("J69").Value = 199860
$ws.Range("L69").Value = 199860
$ws.Range("N69").Value = -201358
$ws.Range("H72").Value = 199860
$ws.Range("J72").Value = 199860
$ws.Range("L72").Value = 599580
$ws.Range("N72").Value = -607068
$ws.Range("H80").Value = 4445.3184
$ws.Range("I80").Value = 3595.889
$ws.Range("K80").Value = 3595.889
$ws.Range("M80").Value = -2597.889
$ws.Range("H83").Value = 4445.3184
$ws.Range("I83").Value = 3595.889
$ws.Range("K83").Value = 17979.445
$ws.Range("M83").Value = -12987.445
$ws.Range("H102").Value = 2969.5652
$ws.Range("I102").Value = 1865.05
$ws.Range("K102").Value = 1865.05
$ws.Range("M102").Value = -243.05
$ws.Range("H113").Value = 3829.3103
$ws.Range("I113").Value = 3168.6191
$ws.Range("J113").Value = 5563.625
$ws.Range("K113").Value = 3168.6191
$ws.Range("L113").Value = 5563.625
$ws.Range("M113").Value = -998.6190999999999
$ws.Range("N113").Value = -9903.625
$ws.Range("H126").Value = 6719.2856
$ws.Range("I126").Value = 6207
$ws.Range("K126").Value = 18621
$ws.Range("M126").Value = -16151

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1301.7
$ws.Range("I22").Value = 1102.4375
$ws.Range("J22").Value = 2098.75
$ws.Range("K22").Value = 1102.4375
$ws.Range("L22").Value = 2098.75
$ws.Range("M22").Value = -807.4375
$ws.Range("N22").Value = -2688.75
$ws.Range("H27").Value = 1301.7
$ws.Range("I27").Value = 1102.4375
$ws.Range("J27").Value = 2098.75
$ws.Range("K27").Value = 1102.4375
$ws.Range("L27").Value = 2098.75
$ws.Range("M27").Value = -995.4375
$ws.Range("N27").Value = -2312.75
$ws.Range("H40").Value = 6464.7646
$ws.Range("I40").Value = 4328.8335
$ws.Range("K40").Value = 4328.8335
$ws.Range("M40").Value = -4192.8335
$ws.Range("H46").Value = 2880.875
$ws.Range("I46").Value = 1605.25
$ws.Range("J46").Value = 4156.5
$ws.Range("K46").Value = 1605.25
$ws.Range("L46").Value = 4156.5
$ws.Range("M46").Value = -1417.25
$ws.Range("N46").Value = -4532.5
$ws.Range("H93").Value = 2206.1
$ws.Range("I93").Value = 1632.125
$ws.Range("K93").Value = 1632.125
$ws.Range("M93").Value = -384.125
$ws.Range("H122").Value = 6166.9287
$ws.Range("I122").Value = 5821.8184
$ws.Range("K122").Value = 17465.4552
$ws.Range("M122").Value = -15015.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2579.0312
$ws.Range("I122").Value = 2294.1333
$ws.Range("J122").Value = 6852.5
$ws.Range("K122").Value = 6882.3999
$ws.Range("L122").Value = 20557.5
$ws.Range("M122").Value = -4432.3999
$ws.Range("N122").Value = -25457.5
